$wb = $excel.ActiveWorkbook

# --- Sheet: Estadisticos 1P ---
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D2").Value = 0
$ws1.Range("E2").Value = 17
$ws1.Range("H2").Value = 6.3

$ws1.Range("D3").Value = 0
$ws1.Range("E3").Value = 8
$ws1.Range("F3").Value = 23
$ws1.Range("G3").Value = 74.19
$ws1.Range("H3").Value = 6.4

$ws1.Range("D4").Value = 0
$ws1.Range("E4").Value = 3
$ws1.Range("H4").Value = 7

$ws1.Range("D5").Value = 0
$ws1.Range("E5").Value = 6
$ws1.Range("F5").Value = 34
$ws1.Range("G5").Value = 85
$ws1.Range("H5").Value = 7.6

$ws1.Range("D6").Value = 0
$ws1.Range("E6").Value = 7
$ws1.Range("H6").Value = 6.4

# --- Sheet: Estadisticos 2P ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("D2").Value = 0
$ws2.Range("E2").Value = 22
$ws2.Range("F2").Value = 14
$ws2.Range("G2").Value = 38.89
$ws2.Range("H2").Value = 6.3

$ws2.Range("D3").Value = 0
$ws2.Range("E3").Value = 10
$ws2.Range("F3").Value = 21
$ws2.Range("G3").Value = 67.73999999999999
$ws2.Range("H3").Value = 6.4

$ws2.Range("D4").Value = 0
$ws2.Range("E4").Value = 7
$ws2.Range("F4").Value = 14
$ws2.Range("G4").Value = 66.67
$ws2.Range("H4").Value = 7

$ws2.Range("D5").Value = 0
$ws2.Range("E5").Value = 10
$ws2.Range("F5").Value = 30
$ws2.Range("G5").Value = 75
$ws2.Range("H5").Value = 7.6

$ws2.Range("D6").Value = 0
$ws2.Range("E6").Value = 10
$ws2.Range("F6").Value = 13
$ws2.Range("G6").Value = 56.52
$ws2.Range("H6").Value = 6.4

# --- Sheet: Estadisticos Final ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("D2").Value = 0
$ws3.Range("E2").Value = 22
$ws3.Range("F2").Value = 14
$ws3.Range("G2").Value = 38.89
$ws3.Range("H2").Value = 6.1

$ws3.Range("D3").Value = 0
$ws3.Range("E3").Value = 10
$ws3.Range("F3").Value = 21
$ws3.Range("G3").Value = 67.73999999999999
$ws3.Range("H3").Value = 6.9

$ws3.Range("D4").Value = 0
$ws3.Range("E4").Value = 7
$ws3.Range("F4").Value = 14
$ws3.Range("G4").Value = 66.67
$ws3.Range("H4").Value = 7.1

$ws3.Range("D5").Value = 0
$ws3.Range("E5").Value = 10
$ws3.Range("F5").Value = 30
$ws3.Range("G5").Value = 75
$ws3.Range("H5").Value = 7.5

$ws3.Range("D6").Value = 0
$ws3.Range("E6").Value = 10
$ws3.Range("F6").Value = 13
$ws3.Range("G6").Value = 56.52
$ws3.Range("H6").Value = 6.2

# --- Sheet: Rescatables ---
$ws4 = $wb.Worksheets.Item("Rescatables")
$ws4.Range("A2").Value = 20330051920326
$ws4.Range("B2").Value = "LUNA"
$ws4.Range("C2").Value = "MORALES"
$ws4.Range("D2").Value = "JESUS ANTONIO"
$ws4.Range("E2").Value = "QUÍMICA II"
$ws4.Range("F2").Value = "2ASV"
$ws4.Range("G2").Value = 2
